$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73, shifting existing rows 73-133 down to 74-134.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new price observation.
$ws.Range("A73").Value = 10
$ws.Range("B73").Value = "Vega Modelo de Temuco"
$ws.Range("C73").Value = "La Araucanía"
$ws.Range("D73").Value = 45264
$ws.Range("E73").Value = 9
$ws.Range("F73").Value = 100112022
$ws.Range("G73").Value = "Arveja Verde"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 50
$ws.Range("K73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("M73").Value = 25000
$ws.Range("N73").Value = '$/saco 25 kilos'
$ws.Range("O73").Value = "Región del Maule"
$ws.Range("P73").Value = 1000
$ws.Range("Q73").Value = 25
$ws.Range("R73").Value = "Hortaliza"
